$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume(1h) columns to stay text so Excel does not
# silently coerce numeric-looking strings (e.g. "1.001") into floating
# point numbers and mangle their original formatting/precision.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "25.866.62"
$ws.Range("E2").Value = "  -0.17%  "

$ws.Range("D3").Value = "1.635.87"
$ws.Range("E3").Value = "  -0.03%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.78%  "

$ws.Range("D5").Value = "213.78"
$ws.Range("E5").Value = "  -0.44%  "

$ws.Range("D6").Value = "0.5009"
$ws.Range("E6").Value = "  -0.46%  "

$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.45%  "

$ws.Range("D8").Value = "0.2557"
$ws.Range("E8").Value = "  -0.56%  "

$ws.Range("D9").Value = "0.06364"
$ws.Range("E9").Value = "  -1.18%  "

$ws.Range("E10").Value = "  +0.08%  "

$ws.Range("D11").Value = "0.07770"
$ws.Range("E11").Value = "  +0.62%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.655.44"
$ws.Range("E12").Value = "  +0.95%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.251"
$ws.Range("E13").Value = "  +0.14%  "

$ws.Range("D14").Value = "1.863.78"
$ws.Range("E14").Value = "  +0.01%  "

$ws.Range("D15").Value = "0.5403"
$ws.Range("E15").Value = "  -0.80%  "

$ws.Range("D16").Value = "0.0₅7822"
$ws.Range("E16").Value = "  -1.55%  "

$ws.Range("D17").Value = "64.45"
$ws.Range("E17").Value = "  +1.57%  "

$ws.Range("D18").Value = "25.902.00"
$ws.Range("E18").Value = "  -0.10%  "

$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  -0.46%  "

$ws.Range("D20").Value = "197.34"
$ws.Range("E20").Value = "  -3.47%  "

$ws.Range("D21").Value = "4.363"
$ws.Range("E21").Value = "  +1.45%  "

$ws.Range("D22").Value = "9.883"
$ws.Range("E22").Value = "  -0.84%  "

$ws.Range("D23").Value = "5.950"
$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").Value = "1.004"
$ws.Range("E24").Value = "  -0.17%  "

$ws.Range("D25").Value = "1.865"
$ws.Range("E25").Value = "  -1.53%  "

$ws.Range("D26").Value = "139.58"
$ws.Range("E26").Value = "  -1.20%  "

$ws.Range("D27").Value = "0.1135"
$ws.Range("E27").Value = "  -1.60%  "

$ws.Range("D28").Value = "6.810"
$ws.Range("E28").Value = "  +1.20%  "

$ws.Range("D29").Value = "15.63"
$ws.Range("E29").Value = "  -1.04%  "

$ws.Range("D30").Value = "1.236"
$ws.Range("E30").Value = "  +0.01%  "

$ws.Range("D31").Value = "0.04873"
$ws.Range("E31").Value = "  -3.61%  "

$ws.Range("D32").Value = "3.243"
$ws.Range("E32").Value = "  -0.45%  "

$ws.Range("D33").Value = "3.174"
$ws.Range("E33").Value = "  -0.34%  "

$ws.Range("E34").Value = "  -0.89%  "

$ws.Range("D35").Value = "2.359"
$ws.Range("E35").Value = "  +1.07%  "

$ws.Range("E36").Value = "  -0.71%  "

$ws.Range("D37").Value = "2.592"
$ws.Range("E37").Value = "  -0.99%  "

$ws.Range("D38").Value = "1.132.24"
$ws.Range("E38").Value = "  -1.72%  "

$ws.Range("D39").Value = "0.5517"
$ws.Range("E39").Value = "  -2.31%  "

$ws.Range("D40").Value = "0.01556"
$ws.Range("E40").Value = "  -0.59%  "

$ws.Range("D41").Value = "1.000"
$ws.Range("E41").Value = "  -0.38%  "

$ws.Range("D42").Value = "5.660"
$ws.Range("E42").Value = "  +0.52%  "

$ws.Range("D43").Value = "0.8129"
$ws.Range("E43").Value = "  -0.23%  "

$ws.Range("D44").Value = "99.22"
$ws.Range("E44").Value = "  -0.26%  "

$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₈123"
$ws.Range("E45").Value = "  +8.30%  "

$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.774.79"
$ws.Range("E46").Value = "  -0.09%  "

$ws.Range("D47").Value = "0.4504"
$ws.Range("E47").Value = "  -0.19%  "

$ws.Range("D48").Value = "1.007"
$ws.Range("E48").Value = "  +0.04%  "

$ws.Range("D49").Value = "55.04"
$ws.Range("E49").Value = "  +0.61%  "

$ws.Range("D50").Value = "0.05073"
$ws.Range("E50").Value = "  +0.82%  "

$ws.Range("E51").Value = "  -0.17%  "
